# Apply "historique" worksheet update:
#  - add a new "Reference" column (G) with values for existing rows 2-6
#  - append new movement rows 7-19 (including column G)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cell for column G (reuse the header style from column A) ---
$ws.Range("G1").Value = "Reference"
$ws.Range("A1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats

# --- New rows data: Date, Produit, Nature, Qte_Mouvement, Qte_Avant, Qte_Apres, Reference ---
$rows = @(
    @("2025-05-26 13:42:18", "Allonge de 600 1E 255282", "Sortie", 3, 192, 189, "2923465824"),
    @("2025-05-26 13:42:26", "Allonge de 600 1E 255282", "Sortie", 10, 189, 179, "2923465824"),
    @("2025-05-26 13:42:31", "Allonge de 600 1E 255282", "Entrée", 10, 179, 189, "2923465824"),
    @("2025-05-26 13:43:08", "Allonge de 600 1E 255282", "Entrée", 14, 189, 203, "2923465824"),
    @("2025-05-26 13:43:16", "Allonge de 600 1E 255282", "Sortie", 56, 203, 147, "2923465824"),
    @("2025-05-27 11:20:30", "Allonge de 600 1E 255282", "Sortie", 1, 147, 146, "2923465824"),
    @("2025-05-27 11:54:38", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Entrée", 6, 1, 7, "7007332946"),
    @("2025-05-27 11:54:44", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Sortie", 1, 7, 6, "7007332946"),
    @("2025-05-27 11:54:46", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Sortie", 1, 6, 5, "7007332946"),
    @("2025-05-27 11:54:57", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Entrée", 1, 5, 6, "7007332946"),
    @("2025-05-27 11:54:59", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Sortie", 1, 6, 5, "7007332946"),
    @("2025-05-27 11:55:01", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Entrée", 1, 5, 6, "7007332946"),
    @("2025-05-27 11:55:03", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Sortie", 1, 6, 5, "7007332946"),
    @("2025-05-27 11:55:03", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Sortie", 1, 5, 4, "7007332946"),
    @("2025-05-27 11:55:06", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Sortie", 1, 4, 3, "7007332946"),
    @("2025-05-27 11:55:08", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Entrée", 1, 3, 4, "7007332946"),
    @("2025-05-27 11:55:12", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Sortie", 1, 4, 3, "7007332946"),
    @("2025-05-27 11:55:21", "POIGNEE EXT REDUITE 5510 BLANC 9016 DTE ", "Entrée", 4, 3, 7, "7007332946")
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}

# --- Fill column G (Reference) as text for rows 2..19 in one pass so the
#     text number-format style is only added/used once. ---
$lastRow = $startRow + $rows.Count - 1
$refRange = $ws.Range("G2:G$lastRow")
$refRange.NumberFormat = "@"
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 7).Value = $rows[$i][6]
}
$refRange.Style = "Normal"
